$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Primary-key columns across every table changed from BIGINT 20 to INT 11
$pkCells = @("A6","A10","A14","A18","A22","A26","A30","A34","A38","A42","A46","A54")
foreach ($cellRef in $pkCells) {
    $ws.Range($cellRef).Value = "INT 11"
}

# Users table: Alias (U) column type changed from VARCHAR 5 to VARCHAR 6
$ws.Range("D6").Value = "VARCHAR 6"

# Update selection to match the saved cursor position
$ws.Range("L6").Select()
